# Daily auto push: insert two new rows of sensor/count data at row 832,
# pushing the existing rows (832-873) down to (834-875).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 832, shifting rows 832:873 down to 834:875.
$ws.Rows.Item(832).Resize(2).Insert()

# Populate the two newly inserted rows with the new data. Column A holds
# the date as plain text (not a real date serial), so force text format
# before writing the value to stop Excel from auto-converting it, then
# restore the default "Normal" style so the cell ends up unstyled, same
# as every other data cell in the sheet.
$ws.Range("A832:A833").NumberFormat = "@"

$ws.Range("A832").Value = "2026/02/16"
$ws.Range("B832").Value = "月"
$ws.Range("C832").Value = 20
$ws.Range("D832").Value = 201

$ws.Range("A833").Value = "2026/02/16"
$ws.Range("B833").Value = "月"
$ws.Range("C833").Value = 22
$ws.Range("D833").Value = 201

$ws.Range("A832:A833").Style = "Normal"
